# The commit swaps the deck's theme-color palette: the single live
# DrawingML theme (serialized to ppt/theme/theme2.xml, the theme
# actually wired to the slide master / presentation) goes from the
# "Integral" palette to the stock "Office Theme" palette.
#
# PowerPoint's automation surface exposes the 12 DrawingML theme slots
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) through
# Slide.ThemeColorScheme.Colors(index).RGB - that collection is backed
# by the one live theme shared by every slide/layout/master, so editing
# it from slide 1 is sufficient.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function RGBVal([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target "Office Theme" color scheme (dk1..folHlink), replacing "Integral".
$officeTheme = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

for ($i = 1; $i -le $officeTheme.Length; $i++) {
    $tcs.Colors($i).RGB = RGBVal($officeTheme[$i - 1])
}
